$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.674.60"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "'3.333.61"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'580.73"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").Value = "'175.76"
$ws.Range("E6").Value = "  -2.01%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("D9").Value = "'3.333.59"
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("D11").Value = "'0.580"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").Value = "'46.27"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").Value = "'0.0000271"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").Value = "'703.02"
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("D15").Value = "'3.883.21"
$ws.Range("E15").Value = "  +1.77%  "
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "'67.762.48"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").Value = "'3.340.78"
$ws.Range("E19").Value = "  +1.47%  "
$ws.Range("D20").Value = "'17.38"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").Value = "'10.99"
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("D22").Value = "'0.894"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").Value = "'5.41"
$ws.Range("E23").Value = "  +4.06%  "
$ws.Range("D24").Value = "'16.95"
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("D25").Value = "'98.30"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("E26").Value = "  -1.77%  "
$ws.Range("E27").Value = "  -2.12%  "
$ws.Range("D28").Value = "'9.44"
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("D29").Value = "'33.23"
$ws.Range("E29").Value = "  +0.29%  "
$ws.Range("D30").Value = "'8.52"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("D31").Value = "'7.12"
$ws.Range("E31").Value = "  +4.91%  "
$ws.Range("D32").Value = "'571.57"
$ws.Range("E32").Value = "  -1.95%  "
$ws.Range("D33").Value = "'10.98"
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'57.32"
$ws.Range("E35").Value = "  +3.52%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").Value = "'3.709.41"
$ws.Range("E37").Value = "  -4.10%  "
$ws.Range("D38").Value = "'3.33"
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("D39").Value = "'34.03"
$ws.Range("E39").Value = "  +5.38%  "
$ws.Range("D40").Value = "'0.130"
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'3.20"
$ws.Range("E41").Value = "  -0.67%  "
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").Value = "'2.65"
$ws.Range("E42").Value = "  +0.52%  "
$ws.Range("D43").Value = "'0.0₃0674"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").Value = "'0.336"
$ws.Range("E44").Value = "  +1.41%  "
$ws.Range("E45").Value = "  -3.05%  "
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("E47").Value = "  +6.29%  "
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("D50").Value = "'1.33"
$ws.Range("E50").Value = "  -5.19%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'128.55"
$ws.Range("E51").Value = "  -0.45%  "
